$wb = $excel.ActiveWorkbook

$newFileName = "0a0e547e-49be-4dea-9408-234e61a303df.md"

function Update-FileNameHyperlink($ws) {
    $ws.Range("A2").Value = $newFileName
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq '$A$2') {
            $hl.TextToDisplay = $newFileName
        }
    }
}

# --- Overview sheet: file-name rename + rolled-up status for both languages ---
$overview = $wb.Worksheets.Item("Overview")
Update-FileNameHyperlink $overview
$overview.Range("B2").Value = "Handoff transform failed"
$overview.Range("C2").Value = "Handoff transform failed"

# --- zh-cn / de-de sheets: rename + status/handoff updates ---
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)

    Update-FileNameHyperlink $ws

    $ws.Range("B2").Value = "Handoff transform failed"

    $toDelete = @()
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq '$C$2') {
            $toDelete += $hl
        }
    }
    foreach ($hl in $toDelete) {
        $hl.Delete()
    }
    $ws.Range("C2").Clear()

    $ws.Range("D2").Value = "0001-01-01 00:00:00"
    $ws.Range("H2").Value = "Ignored"
}
